$d = $word.ActiveDocument

$pairs = @(
    @("2025-06-26 Thursday", "2025-06-27 Friday"),
    @("18×77=1386", "64×86=5504"),
    @("31×79=2449", "21×59=1239"),
    @("94×50=4700", "97×97=9409"),
    @("21×22=462", "63×64=4032"),
    @("96×98=9408", "84×47=3948"),
    @("59×18=1062", "52×62=3224"),
    @("25×29=725", "67×76=5092"),
    @("93×39=3627", "59×88=5192"),
    @("40×46=1840", "34×54=1836"),
    @("59×67=3953", "19×29=551"),
    @("83×99=8217", "78×26=2028"),
    @("82×52=4264", "42×53=2226"),
    @("20×13=260", "83×94=7802"),
    @("67×19=1273", "78×26=2028"),
    @("38×96=3648", "49×95=4655"),
    @("79×88=6952", "53×18=954"),
    @("23×62=1426", "78×47=3666"),
    @("94×70=6580", "11×16=176"),
    @("58×23=1334", "25×99=2475"),
    @("62×12=744", "45×89=4005"),
    @("53×49=2597", "36×33=1188"),
    @("26×58=1508", "60×42=2520"),
    @("38×45=1710", "20×78=1560"),
    @("70×53=3710", "98×50=4900"),
    @("13×12=156", "86×13=1118")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
